# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 21:52"

# --- Update country statistics (rows keep same meaning except 160/161 swap) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 759134
$ws.Range("C4").Value = 20342
$ws.Range("D4").Value = 69927
$ws.Range("E4").Value = 648942
$ws.Range("G4").Value = 1251
$ws.Range("H4").Value = 40265

# Row 7: Francia
$ws.Range("B7").Value = 152894
$ws.Range("C7").Value = 1101
$ws.Range("E7").Value = 96598

# Row 8: Alemania
$ws.Range("B8").Value = 145184
$ws.Range("C8").Value = 1460
$ws.Range("E8").Value = 52598
$ws.Range("G8").Value = 48
$ws.Range("H8").Value = 4586

# Row 35: Dinamarca
$ws.Range("F35").Value = 84

# Row 94: Costa Rica
$ws.Range("B94").Value = 660
$ws.Range("C94").Value = 5
$ws.Range("D94").Value = 112
$ws.Range("E94").Value = 544

# Row 100: Nigeria
$ws.Range("B100").Value = 541
$ws.Range("E100").Value = 356

# Row 131: Ruanda
$ws.Range("B131").Value = 147
$ws.Range("C131").Value = 3
$ws.Range("D131").Value = 76
$ws.Range("E131").Value = 71

# Rows 160/161: Guinea-Bisau overtakes Libia in the ranking, so the two
# rows swap places (country name + stats move together).
$ws.Range("A160").Value = "Guinea-Bisau"
$ws.Range("B160").Value = 50
$ws.Range("C160").Value = 4
$ws.Range("D160").Value = 3
$ws.Range("E160").Value = 47
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0

$ws.Range("A161").Value = "Libia"
$ws.Range("B161").Value = 49
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 11
$ws.Range("E161").Value = 37
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 1

# Row 182: Belice
$ws.Range("D182").Value = 2
$ws.Range("E182").Value = 14
